$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1, J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, bordered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data columns I and J for rows 2-4
$ws.Cells.Item(2, 9).Value = 8
$ws.Cells.Item(2, 10).Value = 8

$ws.Cells.Item(3, 9).Value = 8
$ws.Cells.Item(3, 10).Value = 8

$ws.Cells.Item(4, 9).Value = 6
$ws.Cells.Item(4, 10).Value = 7
